$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This sheet is being reshaped: a new "Avg_" naming convention for
# H1 plus six new Std_* columns (I..N) are introduced, Obs_Prob moves
# from column I to column O, every data row is re-sorted/updated with
# higher-precision recomputed statistics, and rows swap pairwise
# (2<->3, 6<->7, 10<->11) to match the new sort order.
# ------------------------------------------------------------------

# Copy the bold/centered header style (currently on I1) onto the new
# header cells K1:O1 before we repopulate the header text.
$ws.Range("I1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)
[void]($excel.CutCopyMode = $false)

# Clear old data (old layout only spanned A:J) before rewriting it in
# the new A:O layout.
$ws.Range("A1:J13").ClearContents()

# --- Header row (row 1) ---
$ws.Range("A1").Value = "#_Agents"
$ws.Range("B1").Value = "Coverage"
$ws.Range("C1").Value = "Avg_Total_Rounds"
$ws.Range("D1").Value = "Avg_Expl_Cost"
$ws.Range("E1").Value = "Avg_Expl_Eff"
$ws.Range("F1").Value = "Avg_Round_Time"
$ws.Range("G1").Value = "Avg_Agent_Step_Time"
$ws.Range("H1").Value = "Avg_Experiment_Time"
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"
$ws.Range("O1").Value = "Obs_Prob"

# --- Data rows (rows 2-13) ---
# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 58.474
$ws.Cells.Item(2, 4).Value = 58.474
$ws.Cells.Item(2, 5).Value = 2.94401892
$ws.Cells.Item(2, 6).Value = 0.11105262
$ws.Cells.Item(2, 7).Value = 0.11105262
$ws.Cells.Item(2, 8).Value = 6.373386700000001
$ws.Cells.Item(2, 9).Value = 7.853450028790495
$ws.Cells.Item(2, 10).Value = 7.853450028790495
$ws.Cells.Item(2, 11).Value = 0.4102321740655048
$ws.Cells.Item(2, 12).Value = 0.01794558243114963
$ws.Cells.Item(2, 13).Value = 0.01794558243114963
$ws.Cells.Item(2, 14).Value = 0.5277029664222589
$ws.Cells.Item(2, 15).Value = 0.15

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 89.294
$ws.Cells.Item(3, 4).Value = 89.294
$ws.Cells.Item(3, 5).Value = 1.92376268
$ws.Cells.Item(3, 6).Value = 0.0771154
$ws.Cells.Item(3, 7).Value = 0.0771154
$ws.Cells.Item(3, 8).Value = 6.7946048
$ws.Cells.Item(3, 9).Value = 11.29454281477518
$ws.Cells.Item(3, 10).Value = 11.29454281477518
$ws.Cells.Item(3, 11).Value = 0.2510555711216144
$ws.Cells.Item(3, 12).Value = 0.01139126754987282
$ws.Cells.Item(3, 13).Value = 0.01139126754987282
$ws.Cells.Item(3, 14).Value = 0.6884050764772143
$ws.Cells.Item(3, 15).Value = 0.85

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 27.886
$ws.Cells.Item(4, 4).Value = 55.738
$ws.Cells.Item(4, 5).Value = 3.15679716
$ws.Cells.Item(4, 6).Value = 0.18737478
$ws.Cells.Item(4, 7).Value = 0.09368768
$ws.Cells.Item(4, 8).Value = 2.54180798
$ws.Cells.Item(4, 9).Value = 5.62995531008792
$ws.Cells.Item(4, 10).Value = 11.2718883481814
$ws.Cells.Item(4, 11).Value = 0.6425648943227185
$ws.Cells.Item(4, 12).Value = 0.03708545901958052
$ws.Cells.Item(4, 13).Value = 0.01854280898212283
$ws.Cells.Item(4, 14).Value = 0.3899600541673622
$ws.Cells.Item(4, 15).Value = 0.15

# Row 5
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 47.796
$ws.Cells.Item(5, 4).Value = 93.594
$ws.Cells.Item(5, 5).Value = 1.85109512
$ws.Cells.Item(5, 6).Value = 0.12088106
$ws.Cells.Item(5, 7).Value = 0.0604406
$ws.Cells.Item(5, 8).Value = 2.82203772
$ws.Cells.Item(5, 9).Value = 8.259225970225934
$ws.Cells.Item(5, 10).Value = 14.93476677143082
$ws.Cells.Item(5, 11).Value = 0.2939831618564162
$ws.Cells.Item(5, 12).Value = 0.0225175356088581
$ws.Cells.Item(5, 13).Value = 0.01125874271021546
$ws.Cells.Item(5, 14).Value = 0.3627017832494711
$ws.Cells.Item(5, 15).Value = 0.85

# Row 6
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 14.672
$ws.Cells.Item(6, 4).Value = 58.636
$ws.Cells.Item(6, 5).Value = 3.06949904
$ws.Cells.Item(6, 6).Value = 0.23585114
$ws.Cells.Item(6, 7).Value = 0.05896272
$ws.Cells.Item(6, 8).Value = 0.84061064
$ws.Cells.Item(6, 9).Value = 3.768928979339171
$ws.Cells.Item(6, 10).Value = 15.06909616794044
$ws.Cells.Item(6, 11).Value = 0.7651823889608477
$ws.Cells.Item(6, 12).Value = 0.06273600036649107
$ws.Cells.Item(6, 13).Value = 0.01568396212496941
$ws.Cells.Item(6, 14).Value = 0.2334679305626209
$ws.Cells.Item(6, 15).Value = 0.15

# Row 7
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = 0.99986666
$ws.Cells.Item(7, 3).Value = 25.396
$ws.Cells.Item(7, 4).Value = 93.756
$ws.Cells.Item(7, 5).Value = 1.85729004
$ws.Cells.Item(7, 6).Value = 0.1606103
$ws.Cells.Item(7, 7).Value = 0.04015258
$ws.Cells.Item(7, 8).Value = 0.9932987799999999
$ws.Cells.Item(7, 9).Value = 5.102082564305087
$ws.Cells.Item(7, 10).Value = 16.05066077838684
$ws.Cells.Item(7, 11).Value = 0.3301926621223392
$ws.Cells.Item(7, 12).Value = 0.0357451258853444
$ws.Cells.Item(7, 13).Value = 0.0089363824115136
$ws.Cells.Item(7, 14).Value = 0.1877069815701375
$ws.Cells.Item(7, 15).Value = 0.85

# Row 8
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 9.408
$ws.Cells.Item(8, 4).Value = 56.366
$ws.Cells.Item(8, 5).Value = 3.17942708
$ws.Cells.Item(8, 6).Value = 0.27245646
$ws.Cells.Item(8, 7).Value = 0.0454094
$ws.Cells.Item(8, 8).Value = 0.42007648
$ws.Cells.Item(8, 9).Value = 2.364893537325892
$ws.Cells.Item(8, 10).Value = 14.1924944176699
$ws.Cells.Item(8, 11).Value = 0.7595062307534944
$ws.Cells.Item(8, 12).Value = 0.0821180794003196
$ws.Cells.Item(8, 13).Value = 0.01368627118238403
$ws.Cells.Item(8, 14).Value = 0.1485822027200427
$ws.Cells.Item(8, 15).Value = 0.15

# Row 9
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 17.336
$ws.Cells.Item(9, 4).Value = 88.31
$ws.Cells.Item(9, 5).Value = 1.98433102
$ws.Cells.Item(9, 6).Value = 0.1755793
$ws.Cells.Item(9, 7).Value = 0.02926322
$ws.Cells.Item(9, 8).Value = 0.4915637
$ws.Cells.Item(9, 9).Value = 4.356726538614116
$ws.Cells.Item(9, 10).Value = 16.6618466644168
$ws.Cells.Item(9, 11).Value = 0.3854826833894694
$ws.Cells.Item(9, 12).Value = 0.04469419442097462
$ws.Cells.Item(9, 13).Value = 0.007449003230624868
$ws.Cells.Item(9, 14).Value = 0.1244192119679411
$ws.Cells.Item(9, 15).Value = 0.85

# Row 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 6.93
$ws.Cells.Item(10, 4).Value = 55.316
$ws.Cells.Item(10, 5).Value = 3.2261
$ws.Cells.Item(10, 6).Value = 0.24468378
$ws.Cells.Item(10, 7).Value = 0.03058538
$ws.Cells.Item(10, 8).Value = 0.21189334
$ws.Cells.Item(10, 9).Value = 1.68605621826269
$ws.Cells.Item(10, 10).Value = 13.47364918813347
$ws.Cells.Item(10, 11).Value = 0.740832802649669
$ws.Cells.Item(10, 12).Value = 0.0806531619825731
$ws.Cells.Item(10, 13).Value = 0.01008205525336296
$ws.Cells.Item(10, 14).Value = 0.08795885177658991
$ws.Cells.Item(10, 15).Value = 0.15

# Row 11
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 13.476
$ws.Cells.Item(11, 4).Value = 83.838
$ws.Cells.Item(11, 5).Value = 2.090531
$ws.Cells.Item(11, 6).Value = 0.15165912
$ws.Cells.Item(11, 7).Value = 0.01895734
$ws.Cells.Item(11, 8).Value = 0.24756714
$ws.Cells.Item(11, 9).Value = 3.821863063291135
$ws.Cells.Item(11, 10).Value = 15.98094293445003
$ws.Cells.Item(11, 11).Value = 0.407137869147526
$ws.Cells.Item(11, 12).Value = 0.04116813671015452
$ws.Cells.Item(11, 13).Value = 0.005146037478279352
$ws.Cells.Item(11, 14).Value = 0.07564787741857845
$ws.Cells.Item(11, 15).Value = 0.85

# Row 12
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 5.684
$ws.Cells.Item(12, 4).Value = 56.6
$ws.Cells.Item(12, 5).Value = 3.19069958
$ws.Cells.Item(12, 6).Value = 0.22061486
$ws.Cells.Item(12, 7).Value = 0.02206134
$ws.Cells.Item(12, 8).Value = 0.12886902
$ws.Cells.Item(12, 9).Value = 1.54431662156933
$ws.Cells.Item(12, 10).Value = 15.32454512756361
$ws.Cells.Item(12, 11).Value = 0.8121985443348584
$ws.Cells.Item(12, 12).Value = 0.0773864654435651
$ws.Cells.Item(12, 13).Value = 0.007738704536303244
$ws.Cells.Item(12, 14).Value = 0.06717126594890077
$ws.Cells.Item(12, 15).Value = 0.15

# Row 13
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 10.62
$ws.Cells.Item(13, 4).Value = 75.106
$ws.Cells.Item(13, 5).Value = 2.33056984
$ws.Cells.Item(13, 6).Value = 0.13876628
$ws.Cells.Item(13, 7).Value = 0.01387668
$ws.Cells.Item(13, 8).Value = 0.143899
$ws.Cells.Item(13, 9).Value = 2.953141053889627
$ws.Cells.Item(13, 10).Value = 13.95780802309658
$ws.Cells.Item(13, 11).Value = 0.4479694552131678
$ws.Cells.Item(13, 12).Value = 0.04020351098040946
$ws.Cells.Item(13, 13).Value = 0.004020366609271584
$ws.Cells.Item(13, 14).Value = 0.04769751945042264
$ws.Cells.Item(13, 15).Value = 0.85

[void]($ws.Range("A1").Select())
